$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 8: product name/provider renamed, and costs updated
$ws.Range("A8").Value = "Gatimar"
$ws.Range("B8").Value = "Hill"
$ws.Range("C8").Value = 10.0
$ws.Range("D8").Value = 3.0
$ws.Range("E8").Value = 14.0

# Add new row 9
$ws.Range("A9").Value = "Ponedora Fase 2"
$ws.Range("B9").Value = "Comayma"
$ws.Range("C9").Value = 256.0
$ws.Range("D9").Value = 3.0
$ws.Range("E9").Value = 30.0
